$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4459
$ws.Range("I11").Value = 4459
$ws.Range("K11").Value = 4459
$ws.Range("M11").Value = -4319
$ws.Range("H100").Value = 3303.9412
$ws.Range("I100").Value = 2745.8
$ws.Range("K100").Value = 2745.8
$ws.Range("M100").Value = -2204.8
$ws.Range("H132").Value = 3775730
$ws.Range("I132").Value = 4349971
$ws.Range("J132").Value = 2145.5715
$ws.Range("K132").Value = 13049913
$ws.Range("L132").Value = 6436.7145
$ws.Range("M132").Value = -13047383
$ws.Range("N132").Value = -11496.7145
$ws.Range("H138").Value = 3842.058
$ws.Range("I138").Value = 2737.389
$ws.Range("J138").Value = 4134.4707
$ws.Range("K138").Value = 8212.167000000001
$ws.Range("L138").Value = 12403.4121
$ws.Range("M138").Value = -3072.167000000001
$ws.Range("N138").Value = -22683.4121
$ws.Range("H141").Value = 458840.78
$ws.Range("I141").Value = 1332.591
$ws.Range("J141").Value = 2471876.8
$ws.Range("K141").Value = 3997.773
$ws.Range("L141").Value = 7415630.399999999
$ws.Range("M141").Value = 1182.227
$ws.Range("N141").Value = -7425990.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 35500
$ws.Range("J24").Value = 35500
$ws.Range("L24").Value = 35500
$ws.Range("N24").Value = -36248
$ws.Range("H32").Value = 6421.6943
$ws.Range("I32").Value = 4624.4917
$ws.Range("K32").Value = 4624.4917
$ws.Range("M32").Value = -4337.4917
$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 30000
$ws.Range("N96").Value = -35492
$ws.Range("H100").Value = 35500
$ws.Range("J100").Value = 35500
$ws.Range("L100").Value = 35500
$ws.Range("N100").Value = -37664
$ws.Range("H132").Value = 14086482
$ws.Range("I132").Value = 20001102
$ws.Range("J132").Value = 4051.9524
$ws.Range("K132").Value = 60003306
$ws.Range("L132").Value = 12155.8572
$ws.Range("M132").Value = -60000776
$ws.Range("N132").Value = -17215.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4118.906
$ws.Range("I31").Value = 2905.6216
$ws.Range("J31").Value = 6924.625
$ws.Range("K31").Value = 2905.6216
$ws.Range("L31").Value = 6924.625
$ws.Range("M31").Value = -2610.6216
$ws.Range("N31").Value = -7514.625
$ws.Range("H34").Value = 4118.906
$ws.Range("I34").Value = 2905.6216
$ws.Range("J34").Value = 6924.625
$ws.Range("K34").Value = 2905.6216
$ws.Range("L34").Value = 6924.625
$ws.Range("M34").Value = -2703.6216
$ws.Range("N34").Value = -7328.625
$ws.Range("H64").Value = 32500
$ws.Range("J64").Value = 32500
$ws.Range("L64").Value = 32500
$ws.Range("N64").Value = -32996
$ws.Range("H67").Value = 32500
$ws.Range("J67").Value = 32500
$ws.Range("L67").Value = 32500
$ws.Range("N67").Value = -34216
$ws.Range("H132").Value = 1747.8928
$ws.Range("I132").Value = 1354.4889
$ws.Range("J132").Value = 3357.2727
$ws.Range("K132").Value = 4063.4667
$ws.Range("L132").Value = 10071.8181
$ws.Range("M132").Value = -1533.4667
$ws.Range("N132").Value = -15131.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1743.1818
$ws.Range("I39").Value = 300
$ws.Range("J39").Value = 1887.5
$ws.Range("K39").Value = 900
$ws.Range("L39").Value = 5662.5
$ws.Range("M39").Value = -606
$ws.Range("N39").Value = -6250.5
$ws.Range("H55").Value = 2668
$ws.Range("J55").Value = 4130
$ws.Range("L55").Value = 12390
$ws.Range("N55").Value = -12744
$ws.Range("H87").Value = 12646.154
$ws.Range("J87").Value = 15677.777
$ws.Range("L87").Value = 47033.331
$ws.Range("N87").Value = -49529.331
$ws.Range("H90").Value = 12646.154
$ws.Range("J90").Value = 15677.777
$ws.Range("L90").Value = 141099.993
$ws.Range("N90").Value = -153579.993
$ws.Range("H131").Value = 1215.7059
$ws.Range("J131").Value = 1167.1818
$ws.Range("L131").Value = 3501.5454
$ws.Range("N131").Value = -13581.5454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4670.65
$ws.Range("I70").Value = 4650.3076
$ws.Range("J70").Value = 4708.4287
$ws.Range("K70").Value = 4650.3076
$ws.Range("L70").Value = 4708.4287
$ws.Range("M70").Value = -4380.3076
$ws.Range("N70").Value = -5248.4287
$ws.Range("H73").Value = 4670.65
$ws.Range("I73").Value = 4650.3076
$ws.Range("J73").Value = 4708.4287
$ws.Range("K73").Value = 4650.3076
$ws.Range("L73").Value = 4708.4287
$ws.Range("M73").Value = -3714.3076
$ws.Range("N73").Value = -6580.4287
$ws.Range("H132").Value = 2872.5417
$ws.Range("I132").Value = 2387.7646
$ws.Range("J132").Value = 4049.8572
$ws.Range("K132").Value = 7163.293799999999
$ws.Range("L132").Value = 12149.5716
$ws.Range("M132").Value = -4633.293799999999
$ws.Range("N132").Value = -17209.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 25141.334
$ws.Range("I64").Value = 14136
$ws.Range("J64").Value = 28285.715
$ws.Range("K64").Value = 14136
$ws.Range("L64").Value = 28285.715
$ws.Range("M64").Value = -13911
$ws.Range("N64").Value = -28735.715
$ws.Range("H67").Value = 25141.334
$ws.Range("I67").Value = 14136
$ws.Range("J67").Value = 28285.715
$ws.Range("K67").Value = 14136
$ws.Range("L67").Value = 28285.715
$ws.Range("M67").Value = -13356
$ws.Range("N67").Value = -29845.715
$ws.Range("H93").Value = 2468.6553
$ws.Range("I93").Value = 1977.8695
$ws.Range("J93").Value = 4350
$ws.Range("K93").Value = 1977.8695
$ws.Range("L93").Value = 4350
$ws.Range("M93").Value = -729.8695
$ws.Range("N93").Value = -6846
$ws.Range("H132").Value = 2967.743
$ws.Range("I132").Value = 1923.4762
$ws.Range("J132").Value = 4534.143
$ws.Range("K132").Value = 5770.4286
$ws.Range("L132").Value = 13602.429
$ws.Range("M132").Value = -3240.4286
$ws.Range("N132").Value = -18662.429
$ws.Range("H137").Value = 27250
$ws.Range("J137").Value = 27250
$ws.Range("L137").Value = 27250
$ws.Range("N137").Value = -37450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 22074.5
$ws.Range("J63").Value = 22074.5
$ws.Range("L63").Value = 22074.5
$ws.Range("N63").Value = -23322.5
$ws.Range("H66").Value = 22074.5
$ws.Range("J66").Value = 22074.5
$ws.Range("L66").Value = 66223.5
$ws.Range("N66").Value = -72463.5
$ws.Range("H132").Value = 8633.825999999999
$ws.Range("I132").Value = 841.6863
$ws.Range("J132").Value = 30711.555
$ws.Range("K132").Value = 2525.0589
$ws.Range("L132").Value = 92134.66500000001
$ws.Range("M132").Value = 4.941100000000006
$ws.Range("N132").Value = -97194.66500000001
